$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1) Drop the extra trailing columns (X:AG) from the header rows.
#    Clearing the cells removes them from the sheet entirely and
#    shrinks the worksheet dimension / row "spans" automatically.
# ------------------------------------------------------------------
$ws.Range("X1:AG2").Clear()

# ------------------------------------------------------------------
# 2) The Miller-index / pair columns in the header row (row 2) are
#    reordered - rewrite C2:W2 with the new column labels.
# ------------------------------------------------------------------
$ws.Range("C2").Value = "[5, 1, 1]"
$ws.Range("D2").Value = "[4, 2, 2]"
$ws.Range("E2").Value = "[3, 3, 1]"
$ws.Range("F2").Value = "[3, 1, 1]"
$ws.Range("G2").Value = "[1, 1, 1]"
$ws.Range("H2").Value = "[2, 2, 2]"
$ws.Range("I2").Value = "[3, 3, 3]"
$ws.Range("J2").Value = "[2, 0, 0]"
$ws.Range("K2").Value = "[2, 2, 0]"
$ws.Range("L2").Value = "[4, 2, 0]"
$ws.Range("M2").Value = "[4, 0, 0]"
$ws.Range("N2").Value = "1Pair-A"
$ws.Range("O2").Value = "1Pair-B"
$ws.Range("P2").Value = "2Pairs-A"
$ws.Range("Q2").Value = "2Pairs-B"
$ws.Range("R2").Value = "3Pairs-A"
$ws.Range("S2").Value = "3Pairs-B"
$ws.Range("T2").Value = "3Pairs-C"
$ws.Range("U2").Value = "4Pairs"
$ws.Range("V2").Value = "5A4F"
$ws.Range("W2").Value = "MaxUnique"

# ------------------------------------------------------------------
# 3) The scheme names in column B (rows 3:19) come from a shuffled
#    lookup list too - rewrite the labels to the new scheme names.
# ------------------------------------------------------------------
$ws.Range("B3").Value = "Spiral5"
$ws.Range("B4").Value = "RotRing OmegaMax-90"
$ws.Range("B5").Value = "Equal Angle"
$ws.Range("B6").Value = "Tilt Rotate"
$ws.Range("B7").Value = "CLR"
$ws.Range("B8").Value = "Rizzie Hex"
$ws.Range("B9").Value = "Thomas Hex"
$ws.Range("B10").Value = "Tilt Rotate_Partial"
$ws.Range("B11").Value = "RotRing OmegaMax-60"
$ws.Range("B12").Value = "Equal Angle_Partial"
$ws.Range("B13").Value = "Rizzie Hex_Partial"
$ws.Range("B14").Value = "ND Single"
$ws.Range("B15").Value = "RD Single"
$ws.Range("B16").Value = "TD Single"
$ws.Range("B17").Value = "Morris Single"
$ws.Range("B18").Value = "Ring Perpendicular to ND"
$ws.Range("B19").Value = "Ring Perpendicular to RD"

# ------------------------------------------------------------------
# 4) Append the new simulation rows (20:29). Copy the formatting of
#    the last existing row and then overwrite the cell values.
# ------------------------------------------------------------------
$ws.Range("A19:W19").Copy($ws.Range("A20:W20"))
$ws.Range("A20").Value = 18
$ws.Range("B20").Value = "Ring Perpendicular to TD"
$ws.Range("A19:W19").Copy($ws.Range("A21:W21"))
$ws.Range("A21").Value = 19
$ws.Range("B21").Value = "OffsetFTD"
$ws.Range("A19:W19").Copy($ws.Range("A22:W22"))
$ws.Range("A22").Value = 20
$ws.Range("B22").Value = "OffsetATD"
$ws.Range("A19:W19").Copy($ws.Range("A23:W23"))
$ws.Range("A23").Value = 21
$ws.Range("B23").Value = "OffsetF45"
$ws.Range("A19:W19").Copy($ws.Range("A24:W24"))
$ws.Range("A24").Value = 22
$ws.Range("B24").Value = "OffsetA45"
$ws.Range("A19:W19").Copy($ws.Range("A25:W25"))
$ws.Range("A25").Value = 23
$ws.Range("B25").Value = "OffsetFRD"
$ws.Range("A19:W19").Copy($ws.Range("A26:W26"))
$ws.Range("A26").Value = 24
$ws.Range("B26").Value = "OffsetARD"
$ws.Range("A19:W19").Copy($ws.Range("A27:W27"))
$ws.Range("A27").Value = 25
$ws.Range("B27").Value = "Gaussian Quadrature"
$ws.Range("A19:W19").Copy($ws.Range("A28:W28"))
$ws.Range("A28").Value = 26
$ws.Range("B28").Value = "Michael-CCHex"
$ws.Range("A19:W19").Copy($ws.Range("A29:W29"))
$ws.Range("A29").Value = 27
$ws.Range("B29").Value = "Michael-SNHex"

Write-Output "edit complete"
